# Auto-generated edit script: updates leve profit/price calculation columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) across all
# 8 job sheets, reflecting refreshed market-board data from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 314.83334
$ws.Range("I28").Value = 269.07693
$ws.Range("K28").Value = 269.07693
$ws.Range("M28").Value = 215.92307
$ws.Range("H39").Value = 279
$ws.Range("I39").Value = 83
$ws.Range("J39").Value = 421.54544
$ws.Range("K39").Value = 249
$ws.Range("L39").Value = 1264.63632
$ws.Range("M39").Value = 47
$ws.Range("N39").Value = -1856.63632
$ws.Range("H51").Value = 2663.3333
$ws.Range("J51").Value = 1995
$ws.Range("L51").Value = 1995
$ws.Range("N51").Value = -2963
$ws.Range("H62").Value = 3948
$ws.Range("I62").Value = 2627.6
$ws.Range("J62").Value = 6148.6665
$ws.Range("K62").Value = 2627.6
$ws.Range("L62").Value = 6148.6665
$ws.Range("M62").Value = -2003.6
$ws.Range("N62").Value = -7396.6665
$ws.Range("H65").Value = 3948
$ws.Range("I65").Value = 2627.6
$ws.Range("J65").Value = 6148.6665
$ws.Range("K65").Value = 13138
$ws.Range("L65").Value = 30743.3325
$ws.Range("M65").Value = -10018
$ws.Range("N65").Value = -36983.3325
$ws.Range("H98").Value = 870.3
$ws.Range("I98").Value = 709.95654
$ws.Range("K98").Value = 709.95654
$ws.Range("M98").Value = 788.04346
$ws.Range("H112").Value = 1047.4082
$ws.Range("J112").Value = 1057.9362
$ws.Range("L112").Value = 3173.8086
$ws.Range("N112").Value = -5389.8086
$ws.Range("H113").Value = 166669660
$ws.Range("I113").Value = 250001250
$ws.Range("K113").Value = 250001250
$ws.Range("M113").Value = -249997996
$ws.Range("H122").Value = 870.3
$ws.Range("I122").Value = 709.95654
$ws.Range("K122").Value = 2129.86962
$ws.Range("M122").Value = 320.1303800000001
$ws.Range("H135").Value = 21741042
$ws.Range("I135").Value = 1826.0769
$ws.Range("J135").Value = 50002020
$ws.Range("K135").Value = 16434.6921
$ws.Range("L135").Value = 450018180
$ws.Range("M135").Value = -13899.6921
$ws.Range("N135").Value = -450023250
$ws.Range("H137").Value = 1198.2258
$ws.Range("I137").Value = 1137.0454
$ws.Range("J137").Value = 1347.7778
$ws.Range("K137").Value = 3411.1362
$ws.Range("L137").Value = 4043.3334
$ws.Range("M137").Value = -861.1361999999999
$ws.Range("N137").Value = -9143.3334

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8121.5
$ws.Range("I32").Value = 6740.66
$ws.Range("K32").Value = 6740.66
$ws.Range("M32").Value = -6453.66
$ws.Range("H61").Value = 1654.3143
$ws.Range("I61").Value = 1333.0741
$ws.Range("K61").Value = 1333.0741
$ws.Range("M61").Value = -1121.0741
$ws.Range("H63").Value = 5210525
$ws.Range("I63").Value = 2630
$ws.Range("K63").Value = 2630
$ws.Range("M63").Value = -1944
$ws.Range("H66").Value = 5210525
$ws.Range("I66").Value = 2630
$ws.Range("K66").Value = 13150
$ws.Range("M66").Value = -9718
$ws.Range("H132").Value = 11463.712
$ws.Range("I132").Value = 1747.3529
$ws.Range("J132").Value = 29816.834
$ws.Range("K132").Value = 5242.0587
$ws.Range("L132").Value = 89450.50199999999
$ws.Range("M132").Value = -2712.0587
$ws.Range("N132").Value = -94510.50199999999
$ws.Range("H136").Value = 1654.3143
$ws.Range("I136").Value = 1333.0741
$ws.Range("K136").Value = 3999.2223
$ws.Range("M136").Value = -1449.2223

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1844.1072
$ws.Range("I86").Value = 1694.7222
$ws.Range("J86").Value = 2113
$ws.Range("K86").Value = 1694.7222
$ws.Range("L86").Value = 2113
$ws.Range("M86").Value = -571.7221999999999
$ws.Range("N86").Value = -4359
$ws.Range("H89").Value = 1844.1072
$ws.Range("I89").Value = 1694.7222
$ws.Range("J89").Value = 2113
$ws.Range("K89").Value = 8473.610999999999
$ws.Range("L89").Value = 10565
$ws.Range("M89").Value = -2857.610999999999
$ws.Range("N89").Value = -21797
$ws.Range("H107").Value = 1687.2
$ws.Range("I107").Value = 1763.75
$ws.Range("K107").Value = 1763.75
$ws.Range("M107").Value = 156.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1034.1818
$ws.Range("I16").Value = 1123.5555
$ws.Range("J16").Value = 632
$ws.Range("K16").Value = 1123.5555
$ws.Range("L16").Value = 632
$ws.Range("M16").Value = -836.5554999999999
$ws.Range("N16").Value = -1206
$ws.Range("H58").Value = 31269.646
$ws.Range("I58").Value = 1877.8182
$ws.Range("K58").Value = 1877.8182
$ws.Range("M58").Value = -1674.8182
$ws.Range("H113").Value = 1034.1818
$ws.Range("I113").Value = 1123.5555
$ws.Range("J113").Value = 632
$ws.Range("K113").Value = 1123.5555
$ws.Range("L113").Value = 632
$ws.Range("M113").Value = 1046.4445
$ws.Range("N113").Value = -4972
$ws.Range("H132").Value = 1909.6938
$ws.Range("I132").Value = 1506.0526
$ws.Range("J132").Value = 3304.0908
$ws.Range("K132").Value = 4518.1578
$ws.Range("L132").Value = 9912.2724
$ws.Range("M132").Value = -1988.1578
$ws.Range("N132").Value = -14972.2724
$ws.Range("H134").Value = 782.6667
$ws.Range("I134").Value = 649.17645
$ws.Range("J134").Value = 1350
$ws.Range("K134").Value = 1947.52935
$ws.Range("L134").Value = 4050
$ws.Range("M134").Value = 587.4706499999998
$ws.Range("N134").Value = -9120
$ws.Range("H136").Value = 31269.646
$ws.Range("I136").Value = 1877.8182
$ws.Range("K136").Value = 5633.4546
$ws.Range("M136").Value = -3083.4546

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 633.12
$ws.Range("J122").Value = 713.3333
$ws.Range("L122").Value = 6419.9997
$ws.Range("N122").Value = -11319.9997
$ws.Range("H131").Value = 694.5463999999999
$ws.Range("I131").Value = 411.1111
$ws.Range("J131").Value = 723.5341
$ws.Range("K131").Value = 1233.3333
$ws.Range("L131").Value = 2170.6023
$ws.Range("M131").Value = 3806.6667
$ws.Range("N131").Value = -12250.6023

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9745.083000000001
$ws.Range("I113").Value = 11789
$ws.Range("K113").Value = 11789
$ws.Range("M113").Value = -9619
$ws.Range("H122").Value = 3314.25
$ws.Range("I122").Value = 1585.6666
$ws.Range("K122").Value = 4756.9998
$ws.Range("M122").Value = -2306.9998
$ws.Range("H132").Value = 14417.432
$ws.Range("I132").Value = 2940.25
$ws.Range("J132").Value = 45023.25
$ws.Range("K132").Value = 8820.75
$ws.Range("L132").Value = 135069.75
$ws.Range("M132").Value = -6290.75
$ws.Range("N132").Value = -140129.75
$ws.Range("H139").Value = 27250.285
$ws.Range("J139").Value = 27250.285
$ws.Range("L139").Value = 27250.285
$ws.Range("N139").Value = -37530.285

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5222.385
$ws.Range("I61").Value = 1899
$ws.Range("J61").Value = 10539.8
$ws.Range("K61").Value = 1899
$ws.Range("L61").Value = 10539.8
$ws.Range("M61").Value = -1697
$ws.Range("N61").Value = -10943.8
$ws.Range("H113").Value = 5222.385
$ws.Range("I113").Value = 1899
$ws.Range("J113").Value = 10539.8
$ws.Range("K113").Value = 1899
$ws.Range("L113").Value = 10539.8
$ws.Range("M113").Value = 271
$ws.Range("N113").Value = -14879.8
$ws.Range("H122").Value = 895038.25
$ws.Range("I122").Value = 1156561.2
$ws.Range("K122").Value = 3469683.6
$ws.Range("M122").Value = -3467233.6
$ws.Range("H132").Value = 484181.97
$ws.Range("I132").Value = 635323.75
$ws.Range("J132").Value = 5566.3335
$ws.Range("K132").Value = 1905971.25
$ws.Range("L132").Value = 16699.0005
$ws.Range("M132").Value = -1903441.25
$ws.Range("N132").Value = -21759.0005
$ws.Range("H136").Value = 1579.92
$ws.Range("I136").Value = 1395.75
$ws.Range("K136").Value = 4187.25
$ws.Range("M136").Value = -1637.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 961.8333
$ws.Range("I113").Value = 961.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2885.4999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -715.4998999999998
$ws.Range("H119").Value = 49000
$ws.Range("J119").Value = 49000
$ws.Range("L119").Value = 49000
$ws.Range("N113").ClearContents()
$ws.Range("N119").Value = -58676
